$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, B, C, D)
$data = @(
    @(44313, 1, 4, 175.1313485113835),
    @(44314, 0, 2, 87.56567425569177),
    @(44315, 0, 1, 43.78283712784589),
    @(44316, 0, 1, 43.78283712784589),
    @(44317, 0, 1, 43.78283712784589),
    @(44318, 1, 2, 87.56567425569177)
)

$startRow = 239
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy style from the row above for column A (date style)
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
